# Fix build script optimization
# Append a new data row (row 99) to each of the 4 worksheets, mirroring the
# existing "row 98" record structure (time / length / ID / actual-length /
# checksum + their decimal counterparts).

$wb = $excel.ActiveWorkbook

$newRow = 99
$prevRow = 98

# Per-sheet values for the new row, in column order A..I.
$sheetData = @{
    "FE_LFT_#1" = @{
        A = 45885.49461805556
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x00"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 256
        I = 15
    }
    "FE_LFT_#2" = @{
        A = 45885.49461805556
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x10"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 272
        I = 14
    }
    "FE_PLT_#1" = @{
        A = 45885.49461805556
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x5C"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 92
        I = 3
    }
    "FE_PLT_#2" = @{
        A = 45885.49461805556
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x5A"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 90
        I = 3
    }
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetData[$sheetName]

    # Column A carries the date/time number format used by the column above.
    $ws.Cells.Item($newRow, 1).Value = $data.A
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
